# The Pearson logo inline pictures (docPr/cNvPr id="2" and id="4") get
# renamed from image1.png -> image2.png, and the BTec logo inline
# pictures (docPr/cNvPr id="1" and id="3") get renamed from
# image2.jpg -> image1.jpg. These live in the document's headers and
# footers. There is no Name property exposed on InlineShape in the
# Word object model, so we round-trip the flattened WordOpenXML,
# targeting each <wp:docPr>/<pic:cNvPr> element precisely by its
# surrounding id/descr so only the intended elements are touched.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# Pearson Edexcel logo, docPr id="2" (footer1.xml)
$xml = $xml.Replace(
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"/>',
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"/>')

# Pearson Edexcel logo, docPr id="4" (footer2.xml)
$xml = $xml.Replace(
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="4" name="image1.png"/>',
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="4" name="image2.png"/>')

# Pearson Edexcel logo, matching cNvPr (id="0") elements -- both occurrences
$xml = $xml.Replace(
    '<pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>',
    '<pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"/>')

# BTec logo, docPr id="1" (header1.xml)
$xml = $xml.Replace(
    '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"/>',
    '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image1.jpg"/>')

# BTec logo, docPr id="3" (header2.xml)
$xml = $xml.Replace(
    '<wp:docPr descr="BTec_Logo-Orange" id="3" name="image2.jpg"/>',
    '<wp:docPr descr="BTec_Logo-Orange" id="3" name="image1.jpg"/>')

# BTec logo, matching cNvPr (id="0") elements -- both occurrences
$xml = $xml.Replace(
    '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/>',
    '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"/>')

$d.WordOpenXML = $xml

Write-Output "done"
